$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Breed" query used in cell B2 had a variable/query error: it returned an
# extra `Cohort` column (coalesce(co.cohort_description, '') AS `Cohort`) that
# doesn't belong in this particular query. Replace the cell's text with the
# corrected Cypher query (Cohort line removed) per
# "Fixed variables and query errors in Bread from TC01 to TC30".
$fixedBreedQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Bullmastiff']`nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $fixedBreedQuery

# The row heights were re-wrapped/re-flowed after the text edit (one fewer
# line in B2, plus a minor font-metric re-measure of the whole sheet).
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# Re-point the view: author re-saved with the selection on B2 and scrolled
# back to the top-left of the sheet (no more topLeftCell="A4" scroll offset).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
